# fix: number type in output excel
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet Sheet1 -> Evaluasi
$ws.Name = "Evaluasi"

# 2. Header row (A1:J1) loses its bold/bordered/centered look -> back to plain default style
$ws.Range("A1:J1").Style = "Normal"

# 3. Update the evaluation run's results in row 2
$ws.Range("A2").Value = 45937
$ws.Range("B2").Value = 0.008217311197419652
$ws.Range("C2").Value = 1.465681849318622
$ws.Range("D2").Value = 0.01215766150392072
$ws.Range("E2").Value = 0.0000675242033152384
$ws.Range("F2").Value = 29
$ws.Range("G2").Value = 4277.950583559414
$ws.Range("H2").Value = 65.40604393754001
$ws.Range("I2").Value = 54.02125274057797
$ws.Range("J2").Value = 0.7915770596779452

# 4. Date cell keeps a date-style format, now DD/MM/YYYY instead of YYYY-MM-DD
$ws.Range("A2").NumberFormat = "DD/MM/YYYY"

# 5. Numeric metric cells get a proper numeric display format (the "number type" fix)
$ws.Range("B2:J2").NumberFormat = "#,##0.00"
